{"js": "// Remove the trailing \"ToDo List\" section (the heading plus its four bullet\n// items: \"transfert learning from existing frameworks\", \"add autotune to\n// other parameters (layers, dropout, \u2026)\", \"CNN\", \"RNN\") that was appended at\n// the end of the document, while leaving the final blank paragraph in place.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"ToDo List\",\n  \"transfert learning from existing frameworks\",\n  \"add autotune to other parameters (layers, dropout, \\u2026)\",\n  \"CNN\",\n  \"RNN\"\n]);\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (targets.has(text)) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"ToDo List\" heading and its bullet items (transfert learning...,\n# add autotune..., CNN, RNN) that were appended at the end of the document.\n$d = $word.ActiveDocument\n\n$startPara = $null\n$endPara = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($txt -eq \"ToDo List\") { $startPara = $d.Paragraphs.Item($i) }\n    if ($txt -eq \"RNN\") { $endPara = $d.Paragraphs.Item($i) }\n}\n\nif ($startPara -ne $null -and $endPara -ne $null) {\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
